# Refresh the cryptocurrency market-data table (price in column D,
# 1h volume change in column E) with freshly scraped values, and fix
# the HuobiToken/ARBITRUM rows, which were swapped: coin name, link and
# price all move from row 36 to row 37 and vice versa.
#
# Column D mixes two kinds of text: values such as "1.658.69" that use
# '.' as a thousands separator (never valid numbers) and values such as
# "1.007" or "218.30" that happen to parse as plain numbers too. Every
# D-column cell in the source file is stored as TEXT, so for the latter
# group Set-TextCell writes the value via a leading apostrophe (exactly
# what typing '1.007 into Excel does) so it is not auto-converted to a
# numeric cell, then restores the 'Normal' cell style so the quote-prefix
# formatting Excel applies doesn't leave a stray style behind; values
# that already round-trip as text are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Numeric-looking text: force text storage, as Excel would
        # otherwise convert it to a number on entry.
        $Cell.Value = "'" + $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

Set-TextCell $ws.Range("D2") '26.196.61'
$ws.Range("E2").Value = '  -4.34%  '
Set-TextCell $ws.Range("D3") '1.658.69'
$ws.Range("E3").Value = '  -3.00%  '
Set-TextCell $ws.Range("D4") '1.007'
$ws.Range("E4").Value = '  +0.33%  '
Set-TextCell $ws.Range("D5") '218.30'
$ws.Range("E5").Value = '  -2.70%  '
Set-TextCell $ws.Range("D6") '0.5176'
$ws.Range("E6").Value = '  -3.04%  '
$ws.Range("E7").Value = '  +0.30%  '
Set-TextCell $ws.Range("D8") '0.2570'
$ws.Range("E8").Value = '  -4.11%  '
Set-TextCell $ws.Range("D9") '0.06416'
$ws.Range("E10").Value = '  -5.15%  '
Set-TextCell $ws.Range("D11") '0.07788'
$ws.Range("E11").Value = '  +2.30%  '
Set-TextCell $ws.Range("D12") '1.663.50'
$ws.Range("E12").Value = '  -2.96%  '
$ws.Range("E13").Value = '  -5.66%  '
Set-TextCell $ws.Range("D14") '1.886.23'
$ws.Range("E14").Value = '  -3.05%  '
Set-TextCell $ws.Range("D15") '0.5532'
$ws.Range("E15").Value = '  -4.30%  '
$ws.Range("E16").Value = '  -1.73%  '
Set-TextCell $ws.Range("D17") '64.36'
$ws.Range("E17").Value = '  -5.15%  '
Set-TextCell $ws.Range("D18") '26.226.59'
$ws.Range("E18").Value = '  -4.14%  '
Set-TextCell $ws.Range("D19") '1.007'
$ws.Range("E19").Value = '  +0.34%  '
Set-TextCell $ws.Range("D20") '210.95'
$ws.Range("E20").Value = '  -2.92%  '
Set-TextCell $ws.Range("D21") '4.383'
$ws.Range("E21").Value = '  -6.27%  '
$ws.Range("E22").Value = '  -3.69%  '
Set-TextCell $ws.Range("D23") '5.905'
$ws.Range("E23").Value = '  -1.17%  '
Set-TextCell $ws.Range("D25") '144.01'
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("E26").Value = '  +1.24%  '
Set-TextCell $ws.Range("D27") '0.1164'
$ws.Range("E27").Value = '  -4.18%  '
Set-TextCell $ws.Range("D28") '6.968'
$ws.Range("E28").Value = '  -4.19%  '
Set-TextCell $ws.Range("D29") '15.77'
$ws.Range("E29").Value = '  -3.19%  '
Set-TextCell $ws.Range("D30") '0.05274'
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("E31").Value = '  -2.73%  '
Set-TextCell $ws.Range("D32") '3.362'
$ws.Range("E32").Value = '  -4.02%  '
Set-TextCell $ws.Range("D33") '3.222'
$ws.Range("E33").Value = '  -6.10%  '
Set-TextCell $ws.Range("D34") '1.574'
$ws.Range("E34").Value = '  -4.36%  '
Set-TextCell $ws.Range("D35") '2.762'
$ws.Range("E35").Value = '  -4.00%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws.Range("D36") '2.364'
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range("D37") '0.9263'
$ws.Range("E37").Value = '  -2.66%  '
$ws.Range("E38").Value = '  -2.74%  '
Set-TextCell $ws.Range("D39") '1.162.02'
$ws.Range("E39").Value = '  +10.94%  '
$ws.Range("E40").Value = '  -2.81%  '
Set-TextCell $ws.Range("D42") '0.8375'
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("E43").Value = '  -3.24%  '
Set-TextCell $ws.Range("D44") '99.90'
$ws.Range("E44").Value = '  -1.22%  '
Set-TextCell $ws.Range("D45") '1.796.61'
$ws.Range("E45").Value = '  -3.06%  '
$ws.Range("E46").Value = '  -3.56%  '
$ws.Range("E47").Value = '  -0.07%  '
Set-TextCell $ws.Range("D48") '56.09'
$ws.Range("E48").Value = '  -3.38%  '
$ws.Range("E49").Value = '  +0.30%  '
Set-TextCell $ws.Range("D50") '7.875'
$ws.Range("E50").Value = '  -2.59%  '
Set-TextCell $ws.Range("D51") '0.05088'
$ws.Range("E51").Value = '  -2.79%  '
